# Apply the "Updated cryptos list" data refresh to the crypto price table.
# For each changed cell we either:
#   - set .Value directly (safe for text such as names, links and the
#     "  +x.xx%  " volume strings, and for price strings containing more
#     than one "." such as "70.987.09" which Excel cannot parse as a number), or
#   - first set NumberFormat to "@" (Text) and then set .Value, for D-column
#     price strings that look like plain numbers (e.g. "1.00", "0.543"), so
#     Excel keeps the exact text instead of collapsing it to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.987.09"
$ws.Range("E2").Value = "  +6.09%  "
# Row 3
$ws.Range("D3").Value = "3.663.61"
$ws.Range("E3").Value = "  +18.02%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "621.04"
$ws.Range("E5").Value = "  +7.78%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.70"
$ws.Range("E6").Value = "  +2.31%  "
# Row 7
$ws.Range("D7").Value = "3.659.57"
$ws.Range("E7").Value = "  +17.98%  "
# Row 8
$ws.Range("E8").Value = "  -0.07%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  +5.69%  "
# Row 10
$ws.Range("E10").Value = "  +8.05%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.68"
$ws.Range("E11").Value = "  +4.97%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.502"
$ws.Range("E12").Value = "  +7.29%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.41"
$ws.Range("E13").Value = "  +11.86%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000254"
$ws.Range("E14").Value = "  +6.42%  "
# Row 15
$ws.Range("D15").Value = "4.276.60"
$ws.Range("E15").Value = "  +18.01%  "
# Row 16
$ws.Range("D16").Value = "71.019.26"
$ws.Range("E16").Value = "  +6.12%  "
# Row 17
$ws.Range("D17").Value = "3.657.89"
$ws.Range("E17").Value = "  +17.76%  "
# Row 18
$ws.Range("E18").Value = "  +2.07%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  +7.12%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "520.67"
$ws.Range("E20").Value = "  +8.47%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.88"
$ws.Range("E21").Value = "  +0.91%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.22"
$ws.Range("E22").Value = "  +18.24%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.742"
$ws.Range("E23").Value = "  +7.59%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.54"
$ws.Range("E24").Value = "  +13.62%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.58"
$ws.Range("E25").Value = "  +5.85%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.52"
$ws.Range("E26").Value = "  +7.67%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.11"
$ws.Range("E27").Value = "  +10.20%  "
# Row 28
$ws.Range("E28").Value = "  -0.03%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.55"
$ws.Range("E29").Value = "  +11.53%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.18"
$ws.Range("E30").Value = "  +3.65%  "
# Row 31
$ws.Range("E31").Value = "  +12.34%  "
# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000111"
$ws.Range("E32").Value = "  +18.05%  "
# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.62"
$ws.Range("E33").Value = "  +12.91%  "
# Row 34
$ws.Range("E34").Value = "  +4.83%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.10%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.13"
$ws.Range("E36").Value = "  +9.59%  "
# Row 37
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.03"
$ws.Range("E37").Value = "  +9.04%  "
# Row 38
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.350"
$ws.Range("E38").Value = "  +12.35%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  +9.96%  "
# Row 40
$ws.Range("E40").Value = "  +7.13%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.24"
$ws.Range("E41").Value = "  +4.46%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.56"
$ws.Range("E42").Value = "  -6.01%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "430.21"
$ws.Range("E43").Value = "  +15.83%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.83"
$ws.Range("E44").Value = "  +6.11%  "
# Row 45
$ws.Range("D45").Value = "3.114.27"
$ws.Range("E45").Value = "  +11.27%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  +4.62%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0371"
$ws.Range("E47").Value = "  +7.96%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.39"
$ws.Range("E48").Value = "  +10.51%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.80"
$ws.Range("E49").Value = "  +2.32%  "
# Row 50
$ws.Range("E50").Value = "  +0.00%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("E51").Value = "  +11.07%  "
